$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format-only source range (style s=6: no fill, border 3) used to reset the
# fill/border of cells that currently carry the light-blue s=8 style.
$ws.Range("K9").Copy()

# --- Column L (命中消息点 / hit-message point) ---
$ws.Range("L7").Value = "'0.767,1.333"
$ws.Range("L8").Value = "0.8,1.4"
$ws.Range("L9").Value = "0.8667,1.333,1.8"
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("L10").Value = "0.8667,1.333,1.8"
$ws.Range("L11").Value = "'0.767,1.333"
$ws.Range("L12").Value = "0.8,1.4"
$ws.Range("L13").Value = "0.8667,1.333,1.8"
$ws.Range("L14").PasteSpecial(-4122)
$ws.Range("L14").Value = "0.8667,1.333,1.8"
$ws.Range("L15").Value = "'0.767,1.333"
$ws.Range("L16").Value = "0.8,1.4"
$ws.Range("L17").Value = "0.8667,1.333,1.8"
$ws.Range("L18").PasteSpecial(-4122)
$ws.Range("L18").Value = "0.8667,1.333,1.8"
$ws.Range("L19").Value = "'0.767,1.333"
$ws.Range("L20").Value = "0.8,1.4"
$ws.Range("L21").Value = "0.8667,1.333,1.8"
$ws.Range("L22").PasteSpecial(-4122)
$ws.Range("L22").Value = "0.8667,1.333,1.8"
$ws.Range("L23").Value = "'0.767,1.333"
$ws.Range("L24").Value = "0.8,1.4"
$ws.Range("L25").Value = "0.8667,1.333,1.8"
$ws.Range("L26").PasteSpecial(-4122)
$ws.Range("L26").Value = "0.8667,1.333,1.8"
$ws.Range("L27").Value = "'0.767,1.333"
$ws.Range("L28").Value = "0.8,1.4"
$ws.Range("L29").Value = "0.8667,1.333,1.8"
$ws.Range("L30").PasteSpecial(-4122)
$ws.Range("L30").Value = "0.8667,1.333,1.8"
$ws.Range("L31").Value = "'0.767,1.333"
$ws.Range("L32").Value = "0.8,1.4"
$ws.Range("L33").Value = "0.8667,1.333,1.8"
$ws.Range("L34").PasteSpecial(-4122)
$ws.Range("L34").Value = "0.8667,1.333,1.8"

# --- Column M (被击动作编号 / hit-action id) ---
$ws.Range("M7").Value = "bk1,bk2"
$ws.Range("M8").Value = "bk1,bk3"
$ws.Range("M9").Value = "bk1,bk4,bk5"
$ws.Range("M10").PasteSpecial(-4122)
$ws.Range("M10").Value = "bk1,bk3,bk5"
$ws.Range("M11").Value = "bk1,bk2"
$ws.Range("M12").Value = "bk1,bk3"
$ws.Range("M13").Value = "bk1,bk4,bk5"
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("M14").Value = "bk1,bk3,bk5"
$ws.Range("M15").Value = "bk1,bk2"
$ws.Range("M16").Value = "bk1,bk3"
$ws.Range("M17").Value = "bk1,bk4,bk5"
$ws.Range("M18").PasteSpecial(-4122)
$ws.Range("M18").Value = "bk1,bk3,bk5"
$ws.Range("M19").Value = "bk1,bk2"
$ws.Range("M20").Value = "bk1,bk3"
$ws.Range("M21").Value = "bk1,bk4,bk5"
$ws.Range("M22").PasteSpecial(-4122)
$ws.Range("M22").Value = "bk1,bk3,bk5"
$ws.Range("M23").Value = "bk1,bk2"
$ws.Range("M24").Value = "bk1,bk3"
$ws.Range("M25").Value = "bk1,bk4,bk5"
$ws.Range("M26").PasteSpecial(-4122)
$ws.Range("M26").Value = "bk1,bk3,bk5"
$ws.Range("M27").Value = "bk1,bk2"
$ws.Range("M28").Value = "bk1,bk3"
$ws.Range("M29").Value = "bk1,bk4,bk5"
$ws.Range("M30").PasteSpecial(-4122)
$ws.Range("M30").Value = "bk1,bk3,bk5"
$ws.Range("M31").Value = "bk1,bk2"
$ws.Range("M32").Value = "bk1,bk3"
$ws.Range("M33").Value = "bk1,bk4,bk5"
$ws.Range("M34").PasteSpecial(-4122)
$ws.Range("M34").Value = "bk1,bk3,bk5"

$excel.CutCopyMode = $false

# Update the active selection to match the saved view state
$ws.Range("Q15").Select()